$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Coding Phase Defects")
$rng = $ws3.Range("D7")
$rng.Value = "08/03/2025"
Write-Output $rng.Value2
